$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two oldest shipment rows (row 2: PCIU2855920/COSCO ASIA/... and row 3:
# PCIU2514462/APL SALALAH/...) were removed. Deleting row 2 twice removes both
# and shifts everything below up, so the old row 4 (PCIU1438389/KOTA PERWIRA/...)
# becomes the new row 2.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# A new, most-recent shipment row is appended as row 3.
$ws.Range("A3").Value = "PCIU1970956"
$ws.Range("B3").Value = "COSCO ASIA"

# Columns C/D/E hold values that look numeric ("0", "7032011028") but must be
# stored as text (shared strings), matching the rest of the sheet. Assigning
# them directly via .Value would make Excel infer a numeric type, so instead
# write them as a text formula and convert the formula to a static value via
# copy / paste-special-values, which keeps the literal text without tagging
# the cell with any special (e.g. quote-prefix) style.
$ws.Range("C3").Formula = "=""0"""
$ws.Range("C3").Copy()
$ws.Range("C3").PasteSpecial(-4163)

$ws.Range("D3").Formula = "=""7032011028"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)

$ws.Range("E3").Formula = "=""7032011028"""
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$ws.Range("F3").Value = "DEL900009900"
